# Updated data thru Nov25
# Applies the data changes to the "2019" sheet / Table2:
#  - fix a few existing comment cells that had data corrections
#  - append 20 new activity-log rows (706-725) that were recorded after 11/4/2019

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")
$lo = $ws.ListObjects.Item(1)

# Formula text used by the calculated "Z" column of Table2
$zFormula = '=IF(Table2[[#This Row],[Activity]]="Sleep",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,"NA")'

# ---------------------------------------------------------------------------
# 1. First create all 20 new (blank) table rows so every target cell exists.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 20; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# ---------------------------------------------------------------------------
# 2. Prime the shared-string table: write every *new* piece of text in the
#    same order it was first introduced so new strings land on the expected
#    indices.
# ---------------------------------------------------------------------------
$ws.Range("D712").Value = "glucose Drink (75g)"
$ws.Range("D711").Value = "Blue muffins"
$ws.Range("D649").Value = "Fried chicken (389g)"
$ws.Range("D642").Value = "Banana (84g)"
$ws.Range("D641").Value = "Pasta 437g (Pesto 122g)"
$ws.Range("D636").Value = "banana (163g) peanut butter (167g)"
$ws.Range("D713").Value = "Rooibus Red Tea"
$ws.Range("D716").Value = "Chinese fried rice + mabo tofu"
$ws.Range("D717").Value = "Apple-grape juice"
$ws.Range("D721").Value = "apple pie (120g)"
$ws.Range("D722").Value = "Chicken + broccoli + bread"
$ws.Range("D723").Value = "Melotonin (3mg)"
$ws.Range("D725").Value = "Latte (almond milk)"

# ---------------------------------------------------------------------------
# 3. Fix the remaining (non-string-table-order-sensitive) part of the
#    existing rows whose "Comment"/"Start" values were corrected.
# ---------------------------------------------------------------------------
$ws.Range("A641").Value = 43762.559027777781

# ---------------------------------------------------------------------------
# 4. Fill in the rest of each new row (dates/activity/formula), applying
#    formatting copied from a representative existing row first.
# ---------------------------------------------------------------------------

# --- row 706 : Sleep ---
$r = 706
$ws.Range("A700:E700").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("D$r").ClearContents()
$ws.Range("A$r").Value = 43757.910300925927
$ws.Range("B$r").Value = 43758.258333333331
$ws.Range("C$r").Value = "Sleep"
$ws.Range("E$r").Formula = $zFormula

# --- row 707 : Sleep (End cell uses the alternate date style) ---
$r = 707
$ws.Range("A700:E700").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("D$r").ClearContents()
$ws.Range("A$r").Value = 43755.897256944445
$ws.Range("B$r").Value = 43756.313194444447
$ws.Range("B$r").NumberFormat = "m/d/yy h:mm;@"
$ws.Range("C$r").Value = "Sleep"
$ws.Range("E$r").Formula = $zFormula

# --- row 708 : Sleep ---
$r = 708
$ws.Range("A700:E700").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("D$r").ClearContents()
$ws.Range("A$r").Value = 43757.001701388886
$ws.Range("B$r").Value = 43757.234722222223
$ws.Range("C$r").Value = "Sleep"
$ws.Range("E$r").Formula = $zFormula

# --- row 709 : Sleep ---
$r = 709
$ws.Range("A700:E700").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("D$r").ClearContents()
$ws.Range("A$r").Value = 43757.910300925927
$ws.Range("B$r").Value = 43758.258333333331
$ws.Range("C$r").Value = "Sleep"
$ws.Range("E$r").Formula = $zFormula

# --- row 710 : Sleep ---
$r = 710
$ws.Range("A700:E700").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("D$r").ClearContents()
$ws.Range("A$r").Value = 43758.925694444442
$ws.Range("B$r").Value = 43759.231944444444
$ws.Range("C$r").Value = "Sleep"
$ws.Range("E$r").Formula = $zFormula

# --- row 711 : Food (End/Z left blank) ---
$r = 711
$ws.Range("A700:E700").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("E$r").ClearContents()
$ws.Range("A$r").Value = 43758.270833333336
$ws.Range("C$r").Value = "Food"

# --- row 712 : Food (End/Z left blank, alternate date style) ---
$r = 712
$ws.Range("A700:E700").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("E$r").ClearContents()
$ws.Range("B$r").NumberFormat = "m/d/yy h:mm;@"
$ws.Range("A$r").Value = 43757.239583333336
$ws.Range("C$r").Value = "Food"

# --- row 713 : Food (End/Z left blank, alternate date style) ---
$r = 713
$ws.Range("A700:E700").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("E$r").ClearContents()
$ws.Range("B$r").NumberFormat = "m/d/yy h:mm;@"
$ws.Range("A$r").Value = 43759.618055555555
$ws.Range("C$r").Value = "Food"

# --- row 714 : Food ---
$r = 714
$ws.Range("A702:E702").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("B$r").ClearContents()
$ws.Range("A$r").Value = 43759.760416666664
$ws.Range("C$r").Value = "Food"
$ws.Range("D$r").Value = "Rice + chicken curry"
$ws.Range("E$r").Formula = $zFormula

# --- row 715 : Sleep ---
$r = 715
$ws.Range("A700:E700").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("D$r").ClearContents()
$ws.Range("A$r").Value = 43792.942106481481
$ws.Range("B$r").Value = 43793.249305555553
$ws.Range("C$r").Value = "Sleep"
$ws.Range("E$r").Formula = $zFormula

# --- row 716 : Food ---
$r = 716
$ws.Range("A702:E702").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("B$r").ClearContents()
$ws.Range("A$r").Value = 43792.729166666664
$ws.Range("C$r").Value = "Food"
$ws.Range("E$r").Formula = $zFormula

# --- row 717 : Food ---
$r = 717
$ws.Range("A702:E702").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("B$r").ClearContents()
$ws.Range("A$r").Value = 43792.84375
$ws.Range("C$r").Value = "Food"
$ws.Range("E$r").Formula = $zFormula

# --- row 718 : Food ---
$r = 718
$ws.Range("A702:E702").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("B$r").ClearContents()
$ws.Range("A$r").Value = 43793.260416666664
$ws.Range("C$r").Value = "Food"
$ws.Range("D$r").Value = "Latte"
$ws.Range("E$r").Formula = $zFormula

# --- row 719 : Food ---
$r = 719
$ws.Range("A702:E702").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("B$r").ClearContents()
$ws.Range("A$r").Value = 43793.449305555558
$ws.Range("C$r").Value = "Food"
$ws.Range("D$r").Value = "Granola (90g) almond milk + banana"
$ws.Range("E$r").Formula = $zFormula

# --- row 720 : Food ---
$r = 720
$ws.Range("A702:E702").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("B$r").ClearContents()
$ws.Range("A$r").Value = 43793.545138888891
$ws.Range("C$r").Value = "Food"
$ws.Range("D$r").Value = "Chinese fried rice + mabo tofu"
$ws.Range("E$r").Formula = $zFormula

# --- row 721 : Food ---
$r = 721
$ws.Range("A702:E702").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("B$r").ClearContents()
$ws.Range("A$r").Value = 43793.637499999997
$ws.Range("C$r").Value = "Food"
$ws.Range("E$r").Formula = $zFormula

# --- row 722 : Food ---
$r = 722
$ws.Range("A702:E702").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("B$r").ClearContents()
$ws.Range("A$r").Value = 43793.767361111109
$ws.Range("C$r").Value = "Food"
$ws.Range("E$r").Formula = $zFormula

# --- row 723 : Food ---
$r = 723
$ws.Range("A702:E702").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("B$r").ClearContents()
$ws.Range("A$r").Value = 43793.885416666664
$ws.Range("C$r").Value = "Food"
$ws.Range("E$r").Formula = $zFormula

# --- row 724 : Sleep ---
$r = 724
$ws.Range("A700:E700").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("D$r").ClearContents()
$ws.Range("A$r").Value = 43793.917361111111
$ws.Range("B$r").Value = 43794.208333333336
$ws.Range("C$r").Value = "Sleep"
$ws.Range("E$r").Formula = $zFormula

# --- row 725 : Food ---
$r = 725
$ws.Range("A702:E702").Copy()
$ws.Range("A$r:E$r").PasteSpecial(-4122)
$ws.Range("B$r").ClearContents()
$ws.Range("A$r").Value = 43794.21875
$ws.Range("C$r").Value = "Food"
$ws.Range("E$r").Formula = $zFormula

$ws.Range("A726").Select()
